$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rows whose price / volume figures were refreshed (name & link unchanged) ---
# Price column (D) must stay plain text, matching the source data, so force
# the cell format to Text before writing the value (otherwise Excel would
# auto-convert numeric-looking strings like "1.001" into a float).
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "24.207.10"
$ws.Range("E2").Value = "  +14.11%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.672.63"
$ws.Range("E3").Value = "  +8.65%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.001"
$ws.Range("E4").Value = "  +0.01%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "307.55"
$ws.Range("E5").Value = "  +9.27%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.9961"
$ws.Range("E6").Value = "  +3.81%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.3726"
$ws.Range("E7").Value = "  +2.93%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3432"
$ws.Range("E8").Value = "  +8.02%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "48.04"
$ws.Range("E9").Value = "  +17.70%  "

$ws.Range("E10").Value = "  +7.65%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07287"
$ws.Range("E11").Value = "  +7.28%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.9976"
$ws.Range("E12").Value = "  +0.21%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "20.55"
$ws.Range("E13").Value = "  +9.70%  "

$ws.Range("E14").Value = "  +7.50%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "6.760"
$ws.Range("E15").Value = "  +6.66%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.669.61"

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.00001108"
$ws.Range("E17").Value = "  +6.13%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.9966"
$ws.Range("E18").Value = "  +3.80%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06707"
$ws.Range("E19").Value = "  +11.13%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "81.93"
$ws.Range("E20").Value = "  +14.07%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "16.47"
$ws.Range("E21").Value = "  +9.63%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.157"
$ws.Range("E22").Value = "  +8.76%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "12.02"
$ws.Range("E23").Value = "  +5.97%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "24.168.94"
$ws.Range("E24").Value = "  +13.64%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.399"
$ws.Range("E25").Value = "  +4.12%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "3.386"
$ws.Range("E26").Value = "  -8.41%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.680"
$ws.Range("E27").Value = "  +21.61%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "151.99"
$ws.Range("E28").Value = "  +3.31%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "19.54"
$ws.Range("E29").Value = "  +10.50%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.852.17"
$ws.Range("E30").Value = "  +8.94%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "127.63"
$ws.Range("E31").Value = "  +8.51%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "6.337"
$ws.Range("E32").Value = "  +22.61%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.026"
$ws.Range("E33").Value = "  -0.77%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.9875"
$ws.Range("E34").Value = "  +16.35%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.739"
$ws.Range("E35").Value = "  +16.05%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.08426"
$ws.Range("E36").Value = "  +5.27%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "12.43"
$ws.Range("E37").Value = "  +15.98%  "

# --- Rows 38 & 39 swap position (FraxShare <-> InternetComputer(DFINITY)) with refreshed figures ---
$ws.Range("B38").Value = "InternetComputer(DFINITY)"
$ws.Range("C38").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "5.383"
$ws.Range("E38").Value = "  +8.89%  "

$ws.Range("B39").Value = "FraxShare"
$ws.Range("C39").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "8.926"
$ws.Range("E39").Value = "  +15.67%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.06407"
$ws.Range("E40").Value = "  +9.54%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.295"
$ws.Range("E41").Value = "  +6.24%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.02349"
$ws.Range("E42").Value = "  +12.30%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.2117"
$ws.Range("E43").Value = "  +10.76%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.6135"
$ws.Range("E44").Value = "  +12.50%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.9953"
$ws.Range("E45").Value = "  +3.69%  "

# --- Rows 46 & 47 swap position (PancakeSwap <-> EnergySwap) with refreshed figures ---
$ws.Range("B46").Value = "EnergySwap"
$ws.Range("C46").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "13.24"
$ws.Range("E46").Value = "  +6.19%  "

$ws.Range("B47").Value = "PancakeSwap"
$ws.Range("C47").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.803"
$ws.Range("E47").Value = "  +6.96%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.5981"
$ws.Range("E48").Value = "  +9.87%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "127.18"
$ws.Range("E49").Value = "  +4.61%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.029"
$ws.Range("E50").Value = "  +8.65%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.07129"
$ws.Range("E51").Value = "  +7.65%  "
